# Update NATMI LR-pair output (Ccl11-Cxcr3) with recomputed TPM-derived values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.738254
$ws.Range("H2").Value = 2.214762
$ws.Range("I2").Value = 0.005691320045803731
$ws.Range("J2").Value = 0.005691320045803731
$ws.Range("O2").Value = 0.02773017886769741
$ws.Range("P2").Value = 0.02773017886769741
$ws.Range("Q2").Value = 0.03834589709866666
$ws.Range("R2").Value = 0.3451130738879999
$ws.Range("S2").Value = 0.0001578213228634493
$ws.Range("T2").Value = 0.0001578213228634493

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.738254
$ws.Range("H3").Value = 2.214762
$ws.Range("I3").Value = 0.005691320045803731
$ws.Range("J3").Value = 0.005691320045803731
$ws.Range("M3").Value = 1.821156333333333
$ws.Range("N3").Value = 5.463469
$ws.Range("O3").Value = 0.9722698211323025
$ws.Range("P3").Value = 0.9722698211323026
$ws.Range("Q3").Value = 1.344475947708667
$ws.Range("R3").Value = 12.100283529378
$ws.Range("S3").Value = 0.005533498722940282
$ws.Range("T3").Value = 0.005533498722940282

$ws.Range("I4").Value = 0.9440493064670392
$ws.Range("J4").Value = 0.9440493064670391
$ws.Range("O4").Value = 0.02773017886769741
$ws.Range("P4").Value = 0.02773017886769741
$ws.Range("Q4").Value = 6.360636420112001
$ws.Range("R4").Value = 57.245727781008
$ws.Range("S4").Value = 0.02617865612825669
$ws.Range("T4").Value = 0.02617865612825668

$ws.Range("I5").Value = 0.9440493064670392
$ws.Range("J5").Value = 0.9440493064670391
$ws.Range("M5").Value = 1.821156333333333
$ws.Range("N5").Value = 5.463469
$ws.Range("O5").Value = 0.9722698211323025
$ws.Range("P5").Value = 0.9722698211323026
$ws.Range("Q5").Value = 223.015324350247
$ws.Range("R5").Value = 2007.137919152223
$ws.Range("S5").Value = 0.9178706503387825
$ws.Range("T5").Value = 0.9178706503387825

$ws.Range("G6").Value = 5.698467
$ws.Range("H6").Value = 17.095401
$ws.Range("I6").Value = 0.0439304080539368
$ws.Range("J6").Value = 0.04393040805393679
$ws.Range("O6").Value = 0.02773017886769741
$ws.Range("P6").Value = 0.02773017886769741
$ws.Range("Q6").Value = 0.295985973936
$ws.Range("R6").Value = 2.663873765423999
$ws.Range("S6").Value = 0.001218198073066602
$ws.Range("T6").Value = 0.001218198073066602

$ws.Range("G7").Value = 5.698467
$ws.Range("H7").Value = 17.095401
$ws.Range("I7").Value = 0.0439304080539368
$ws.Range("J7").Value = 0.04393040805393679
$ws.Range("M7").Value = 1.821156333333333
$ws.Range("N7").Value = 5.463469
$ws.Range("O7").Value = 0.9722698211323025
$ws.Range("P7").Value = 0.9722698211323026
$ws.Range("Q7").Value = 10.377799267341
$ws.Range("R7").Value = 93.40019340606899
$ws.Range("S7").Value = 0.04271220998087019
$ws.Range("T7").Value = 0.04271220998087019

$ws.Range("G8").Value = 0.8209666666666666
$ws.Range("H8").Value = 2.4629
$ws.Range("I8").Value = 0.006328965433220369
$ws.Range("J8").Value = 0.006328965433220369
$ws.Range("O8").Value = 0.02773017886769741
$ws.Range("P8").Value = 0.02773017886769741
$ws.Range("Q8").Value = 0.04264210328888889
$ws.Range("R8").Value = 0.3837789295999999
$ws.Range("S8").Value = 0.0001755033435106749
$ws.Range("T8").Value = 0.0001755033435106749

$ws.Range("G9").Value = 0.8209666666666666
$ws.Range("H9").Value = 2.4629
$ws.Range("I9").Value = 0.006328965433220369
$ws.Range("J9").Value = 0.006328965433220369
$ws.Range("M9").Value = 1.821156333333333
$ws.Range("N9").Value = 5.463469
$ws.Range("O9").Value = 0.9722698211323025
$ws.Range("P9").Value = 0.9722698211323026
$ws.Range("Q9").Value = 1.495108644455555
$ws.Range("R9").Value = 13.4559778001
$ws.Range("S9").Value = 0.006153462089709694
$ws.Range("T9").Value = 0.006153462089709695

